$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "2025/10/06"
$ws.Range("A68").Style = "Normal"
$ws.Range("B68").Value = "月"
$ws.Range("C68").Value = 5
$ws.Range("D68").Value = 201
